# "Finishing touches for presentation."
#
# 1. Slide 1  - subtitle: append "- Group 1" to the team names.
# 2. Slide 2  - credits block: append "- Group 1" to "DSCI 445 ".
# 3. Slide 4  - data bullet: broaden the data-source description.
# 4. Slide 4  - data bullet: drop the trailing "games" word.
# 5. Slide 7  - image caption (not the title): rename to "Total fights correlation".
# 6. Slide 8  - bullet list: add a new "Cross Validation!" bullet after "BIC".

$p = $ppt.ActivePresentation

# --- Slide 1: Subtitle shape - append " - Group 1" to the team names line ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
[void]($sh1.TextFrame.TextRange.Runs(1).Text = "Chien Lin Jason Nero Sarah Sublett - Group 1")

# --- Slide 2: Content placeholder - "DSCI 445 " paragraph ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
[void]($sh2.TextFrame.TextRange.Paragraphs(7).Runs(1).Text = "DSCI 445 - Group 1 ")

# --- Slide 4: Content placeholder - data source + game count paragraphs ---
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
[void]($sh4.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Utilizing fight data from hockeyfights.com and general hockey statistics from hockeydb.com and hockey-reference.com")
[void]($sh4.TextFrame.TextRange.Paragraphs(6).Runs(1).Text = "1,312 games per regular season, with each team playing 82")

# --- Slide 7: Caption textbox below the zoomed correlation map image (title stays unchanged) ---
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(3)
[void]($sh7.TextFrame.TextRange.Runs(1).Text = "Total fights correlation")

# --- Slide 8: Content placeholder - add "Cross Validation!" bullet after "BIC" ---
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
[void]($sh8.TextFrame.TextRange.InsertAfter("`rCross Validation!"))
